$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2931.158
$ws.Range("I62").Value = 2245.6155
$ws.Range("J62").Value = 4416.5
$ws.Range("K62").Value = 2245.6155
$ws.Range("L62").Value = 4416.5
$ws.Range("M62").Value = -1621.6155
$ws.Range("N62").Value = -5664.5

$ws.Range("H65").Value = 2931.158
$ws.Range("I65").Value = 2245.6155
$ws.Range("J65").Value = 4416.5
$ws.Range("K65").Value = 11228.0775
$ws.Range("L65").Value = 22082.5
$ws.Range("M65").Value = -8108.077499999999
$ws.Range("N65").Value = -28322.5

$ws.Range("H98").Value = 3882.5
$ws.Range("I98").Value = 1645.8334
$ws.Range("J98").Value = 5799.643
$ws.Range("K98").Value = 1645.8334
$ws.Range("L98").Value = 5799.643
$ws.Range("M98").Value = -147.8334
$ws.Range("N98").Value = -8795.643

$ws.Range("H116").Value = 564005.6
$ws.Range("J116").Value = 10817.167
$ws.Range("L116").Value = 10817.167
$ws.Range("N116").Value = -17701.167

$ws.Range("H122").Value = 3882.5
$ws.Range("I122").Value = 1645.8334
$ws.Range("J122").Value = 5799.643
$ws.Range("K122").Value = 4937.5002
$ws.Range("L122").Value = 17398.929
$ws.Range("M122").Value = -2487.5002
$ws.Range("N122").Value = -22298.929

$ws.Range("H131").Value = 3562.9412
$ws.Range("J131").Value = 4287.222
$ws.Range("L131").Value = 12861.666
$ws.Range("N131").Value = -22941.666

$ws.Range("H132").Value = 31255398
$ws.Range("I132").Value = 40005456
$ws.Range("J132").Value = 5185.7144
$ws.Range("K132").Value = 120016368
$ws.Range("L132").Value = 15557.1432
$ws.Range("M132").Value = -120013838
$ws.Range("N132").Value = -20617.1432

$ws.Range("H135").Value = 969.2
$ws.Range("I135").Value = 483.33334
$ws.Range("J135").Value = 2426.8
$ws.Range("K135").Value = 4350.00006
$ws.Range("L135").Value = 21841.2
$ws.Range("M135").Value = -1815.00006
$ws.Range("N135").Value = -26911.2

$ws.Range("H137").Value = 2383993.2
$ws.Range("I137").Value = 4330687.5
$ws.Range("J137").Value = 4700.3335
$ws.Range("K137").Value = 12992062.5
$ws.Range("L137").Value = 14101.0005
$ws.Range("M137").Value = -12989512.5
$ws.Range("N137").Value = -19201.0005

$ws.Range("H138").Value = 2486.71
$ws.Range("I138").Value = 651.7059
$ws.Range("K138").Value = 1955.1177
$ws.Range("M138").Value = 3184.8823

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 5634.952
$ws.Range("I74").Value = 7578.3335
$ws.Range("J74").Value = 3043.7778
$ws.Range("K74").Value = 7578.3335
$ws.Range("L74").Value = 3043.7778
$ws.Range("M74").Value = -6704.3335
$ws.Range("N74").Value = -4791.7778

$ws.Range("H77").Value = 5634.952
$ws.Range("I77").Value = 7578.3335
$ws.Range("J77").Value = 3043.7778
$ws.Range("K77").Value = 37891.6675
$ws.Range("L77").Value = 15218.889
$ws.Range("M77").Value = -33523.6675
$ws.Range("N77").Value = -23954.889

$ws.Range("H122").Value = 2661.9167
$ws.Range("I122").Value = 1494.3
$ws.Range("K122").Value = 4482.9
$ws.Range("M122").Value = -2032.9

$ws.Range("H132").Value = 3344.5
$ws.Range("J132").Value = 5337
$ws.Range("L132").Value = 16011
$ws.Range("N132").Value = -21071

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 2515576.8
$ws.Range("I7").Value = 11995
$ws.Range("J7").Value = 2873231.2
$ws.Range("K7").Value = 11995
$ws.Range("L7").Value = 2873231.2
$ws.Range("M7").Value = -11882
$ws.Range("N7").Value = -2873457.2

$ws.Range("H87").Value = 41800
$ws.Range("J87").Value = 41800
$ws.Range("L87").Value = 41800
$ws.Range("N87").Value = -44296

$ws.Range("H90").Value = 41800
$ws.Range("J90").Value = 41800
$ws.Range("L90").Value = 125400
$ws.Range("N90").Value = -137880

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 3999.8
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 3999.8
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 3999.8
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -4223.8

$ws.Range("H31").Value = 3934.15
$ws.Range("I31").Value = 1267
$ws.Range("J31").Value = 6601.3
$ws.Range("K31").Value = 1267
$ws.Range("L31").Value = 6601.3
$ws.Range("M31").Value = -972
$ws.Range("N31").Value = -7191.3

$ws.Range("H34").Value = 3934.15
$ws.Range("I34").Value = 1267
$ws.Range("J34").Value = 6601.3
$ws.Range("K34").Value = 1267
$ws.Range("L34").Value = 6601.3
$ws.Range("M34").Value = -1065
$ws.Range("N34").Value = -7005.3

$ws.Range("H58").Value = 3152.3333
$ws.Range("I58").Value = 1819.1837
$ws.Range("K58").Value = 1819.1837
$ws.Range("M58").Value = -1616.1837

$ws.Range("H99").Value = 6063549
$ws.Range("I99").Value = 11112539
$ws.Range("J99").Value = 4761.3335
$ws.Range("K99").Value = 11112539
$ws.Range("L99").Value = 4761.3335
$ws.Range("M99").Value = -11111041
$ws.Range("N99").Value = -7757.3335

$ws.Range("H105").Value = 1403.0834
$ws.Range("I105").Value = 1077
$ws.Range("J105").Value = 2544.375
$ws.Range("K105").Value = 1077
$ws.Range("L105").Value = 2544.375
$ws.Range("M105").Value = 670
$ws.Range("N105").Value = -6038.375

$ws.Range("H122").Value = 2199.4285
$ws.Range("I122").Value = 1149.3334
$ws.Range("J122").Value = 8500
$ws.Range("K122").Value = 3448.0002
$ws.Range("L122").Value = 25500
$ws.Range("M122").Value = -998.0001999999999
$ws.Range("N122").Value = -30400

$ws.Range("H126").Value = 6063549
$ws.Range("I126").Value = 11112539
$ws.Range("J126").Value = 4761.3335
$ws.Range("K126").Value = 33337617
$ws.Range("L126").Value = 14284.0005
$ws.Range("M126").Value = -33335147
$ws.Range("N126").Value = -19224.0005

$ws.Range("H134").Value = 6270.609
$ws.Range("I134").Value = 7635.067
$ws.Range("J134").Value = 3712.25
$ws.Range("K134").Value = 22905.201
$ws.Range("L134").Value = 11136.75
$ws.Range("M134").Value = -20370.201
$ws.Range("N134").Value = -16206.75

$ws.Range("H136").Value = 3152.3333
$ws.Range("I136").Value = 1819.1837
$ws.Range("K136").Value = 5457.551100000001
$ws.Range("M136").Value = -2907.551100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 3378932.8
$ws.Range("I113").Value = 555.9524
$ws.Range("J113").Value = 7813052
$ws.Range("K113").Value = 1667.8572
$ws.Range("L113").Value = 23439156
$ws.Range("M113").Value = 502.1428000000001
$ws.Range("N113").Value = -23443496

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 15271.143
$ws.Range("J39").Value = 15271.143
$ws.Range("L39").Value = 15271.143
$ws.Range("N39").Value = -16335.143

$ws.Range("H102").Value = 2420.75
$ws.Range("I102").Value = 1553.8235
$ws.Range("K102").Value = 1553.8235
$ws.Range("M102").Value = 68.17650000000003

$ws.Range("H122").Value = 4965.8
$ws.Range("I122").Value = 1956.3334
$ws.Range("J122").Value = 9480
$ws.Range("K122").Value = 5869.0002
$ws.Range("L122").Value = 28440
$ws.Range("M122").Value = -3419.0002
$ws.Range("N122").Value = -33340

$ws.Range("H126").Value = 3238.08
$ws.Range("I126").Value = 2909.234
$ws.Range("J126").Value = 4339
$ws.Range("K126").Value = 8727.701999999999
$ws.Range("L126").Value = 13017
$ws.Range("M126").Value = -6257.701999999999
$ws.Range("N126").Value = -17957

$ws.Range("H132").Value = 3872.7778
$ws.Range("I132").Value = 1734.5555
$ws.Range("J132").Value = 6011
$ws.Range("K132").Value = 5203.666499999999
$ws.Range("L132").Value = 18033
$ws.Range("M132").Value = -2673.666499999999
$ws.Range("N132").Value = -23093

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3183.5356
$ws.Range("I7").Value = 2568.9333
$ws.Range("J7").Value = 3892.6924
$ws.Range("K7").Value = 2568.9333
$ws.Range("L7").Value = 3892.6924
$ws.Range("M7").Value = -2456.9333
$ws.Range("N7").Value = -4116.6924

$ws.Range("H22").Value = 1478.7142
$ws.Range("I22").Value = 1450.2
$ws.Range("J22").Value = 1550
$ws.Range("K22").Value = 1450.2
$ws.Range("L22").Value = 1550
$ws.Range("M22").Value = -1155.2
$ws.Range("N22").Value = -2140

$ws.Range("H27").Value = 1478.7142
$ws.Range("I27").Value = 1450.2
$ws.Range("J27").Value = 1550
$ws.Range("K27").Value = 1450.2
$ws.Range("L27").Value = 1550
$ws.Range("M27").Value = -1343.2
$ws.Range("N27").Value = -1764

$ws.Range("H40").Value = 4686.364
$ws.Range("I40").Value = 3270.5881
$ws.Range("J40").Value = 9500
$ws.Range("K40").Value = 3270.5881
$ws.Range("L40").Value = 9500
$ws.Range("M40").Value = -3134.5881
$ws.Range("N40").Value = -9772

$ws.Range("H122").Value = 3740.8667
$ws.Range("I122").Value = 1488
$ws.Range("J122").Value = 6315.5713
$ws.Range("K122").Value = 4464
$ws.Range("L122").Value = 18946.7139
$ws.Range("M122").Value = -2014
$ws.Range("N122").Value = -23846.7139

$ws.Range("H126").Value = 3183.5356
$ws.Range("I126").Value = 2568.9333
$ws.Range("J126").Value = 3892.6924
$ws.Range("K126").Value = 7706.7999
$ws.Range("L126").Value = 11678.0772
$ws.Range("M126").Value = -5236.7999
$ws.Range("N126").Value = -16618.0772

$ws.Range("H132").Value = 4313.3335
$ws.Range("I132").Value = 1675.375
$ws.Range("J132").Value = 7328.143
$ws.Range("K132").Value = 5026.125
$ws.Range("L132").Value = 21984.429
$ws.Range("M132").Value = -2496.125
$ws.Range("N132").Value = -27044.429

$ws.Range("H136").Value = 3577.75
$ws.Range("I136").Value = 1209.15
$ws.Range("J136").Value = 7525.4165
$ws.Range("K136").Value = 3627.45
$ws.Range("L136").Value = 22576.2495
$ws.Range("M136").Value = -1077.45
$ws.Range("N136").Value = -27676.2495

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 160272.14
$ws.Range("I62").Value = 4380.4
$ws.Range("J62").Value = 550001.5
$ws.Range("K62").Value = 4380.4
$ws.Range("L62").Value = 550001.5
$ws.Range("M62").Value = -3756.4
$ws.Range("N62").Value = -551249.5

$ws.Range("H65").Value = 160272.14
$ws.Range("I65").Value = 4380.4
$ws.Range("J65").Value = 550001.5
$ws.Range("K65").Value = 21902
$ws.Range("L65").Value = 2750007.5
$ws.Range("M65").Value = -18782
$ws.Range("N65").Value = -2756247.5

$ws.Range("H122").Value = 2767.7837
$ws.Range("I122").Value = 1456
$ws.Range("J122").Value = 5189.5386
$ws.Range("K122").Value = 4368
$ws.Range("L122").Value = 15568.6158
$ws.Range("M122").Value = -1918
$ws.Range("N122").Value = -20468.6158

$ws.Range("H126").Value = 889116.8
$ws.Range("I126").Value = 1829.2858
$ws.Range("J126").Value = 2131319.5
$ws.Range("K126").Value = 5487.857400000001
$ws.Range("L126").Value = 6393958.5
$ws.Range("M126").Value = -3017.857400000001
$ws.Range("N126").Value = -6398898.5

$ws.Range("H132").Value = 7409290.5
$ws.Range("I132").Value = 1427.3235
$ws.Range("J132").Value = 30306322
$ws.Range("K132").Value = 4281.970499999999
$ws.Range("L132").Value = 90918966
$ws.Range("M132").Value = -1751.970499999999
$ws.Range("N132").Value = -90924026
